# 06.DanhSachChucNang.xlsx - apply commit changes
# - Update the sheet view (scroll position / selection)
# - Add "Hoan thanh" (completion) percentage values to E12, E24, E25
# - Move the "Phan cong" (assignment) values from column G to column F on rows 24-25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update sheetView: topLeftCell A31 -> A4, selection E44 -> E15 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
[void]$ws.Range("E15").Select()

# --- Row 12: set completion (E12) to 90% ---
$ws.Range("E12").Value2 = 0.9

# --- Row 24: set completion (E24) to 90%, move assignee from G24 to F24 ---
$ws.Range("E24").Value2 = 0.9
$ws.Range("F24").Value2 = $ws.Range("G24").Value2
$ws.Range("G24").Clear()

# --- Row 25: set completion (E25) to 100%, move assignee from G25 to F25 ---
$ws.Range("E25").Value2 = 1
$ws.Range("F25").Value2 = $ws.Range("G25").Value2
$ws.Range("G25").Clear()
